$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.434.78"
$ws.Range("E2").Value = "  +0.23%  "

$ws.Range("D3").Value = "3.839.42"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'709.51"
$ws.Range("E5").Value = "  +1.17%  "

$ws.Range("D6").Value = "'173.60"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "3.837.86"
$ws.Range("E7").Value = "  +0.78%  "

$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  -0.30%  "

$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("D11").Value = "'7.32"
$ws.Range("E11").Value = "  +0.49%  "

$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").Value = "'0.0000258"
$ws.Range("E13").Value = "  -0.28%  "

$ws.Range("D14").Value = "'37.05"

$ws.Range("D15").Value = "4.486.63"
$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("D16").Value = "3.847.89"
$ws.Range("E16").Value = "  +1.17%  "

$ws.Range("D17").Value = "71.377.01"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").Value = "'7.28"
$ws.Range("E18").Value = "  +0.83%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'17.50"
$ws.Range("E19").Value = "  -1.23%  "

$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "'0.115"
$ws.Range("E20").Value = "  +0.35%  "

$ws.Range("D21").Value = "'498.75"
$ws.Range("E21").Value = "  +3.60%  "

$ws.Range("D22").Value = "'10.77"
$ws.Range("E22").Value = "  -1.57%  "

$ws.Range("D23").Value = "'0.738"
$ws.Range("E23").Value = "  +3.38%  "

$ws.Range("D24").Value = "'85.49"
$ws.Range("E24").Value = "  +1.26%  "

$ws.Range("E25").Value = "  +1.85%  "

$ws.Range("E26").Value = "  +0.76%  "

$ws.Range("D27").Value = "'12.23"
$ws.Range("E27").Value = "  -0.83%  "

$ws.Range("D28").Value = "3.994.92"
$ws.Range("E28").Value = "  +0.90%  "

$ws.Range("D29").Value = "'2.12"
$ws.Range("E29").Value = "  -2.44%  "

$ws.Range("D30").Value = "'3.15"
$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("D32").Value = "'7.54"
$ws.Range("E32").Value = "  -1.22%  "

$ws.Range("E33").Value = "  -2.41%  "

$ws.Range("D34").Value = "'29.54"
$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("E35").Value = "  -3.93%  "

$ws.Range("E36").Value = "  -0.27%  "

$ws.Range("D37").Value = "3.804.19"
$ws.Range("E37").Value = "  +1.20%  "

$ws.Range("D38").Value = "'0.996"
$ws.Range("E38").Value = "  -1.39%  "

$ws.Range("E39").Value = "  +0.23%  "

$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'3.40"
$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'6.05"
$ws.Range("E41").Value = "  +0.27%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'1.04"
$ws.Range("E42").Value = "  +5.15%  "

$ws.Range("D43").Value = "'2.30"
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("E46").Value = "  +0.46%  "

$ws.Range("D47").Value = "'164.26"
$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("D48").Value = "'431.01"
$ws.Range("E48").Value = "  +3.56%  "

$ws.Range("D49").Value = "'49.07"
$ws.Range("E49").Value = "  +0.47%  "

$ws.Range("D50").Value = "'8.78"
$ws.Range("E50").Value = "  +1.60%  "

$ws.Range("E51").Value = "  -0.06%  "

# Remove the temporary text-forcing style picked up from the
# leading-apostrophe assignments above, restoring default styling.
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
